$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 46 - S07_G01_TF001
$ws.Range("F46").Value = "BacktestDetailChart added to Backtests page (price+equity+projection with trade markers)."
$ws.Range("G46").Value = "implemented"

# Row 47 - S07_G01_TF002
$ws.Range("F47").Value = 'Buy/sell markers and synchronised price/equity panes wired to /api/backtests/{id}/chart-data`.'
$ws.Range("G47").Value = "implemented"

# Row 48 - S07_G01_TF003
$ws.Range("F48").Value = "Projection overlay rendered in equity pane based on unrealised what-if equity path."
$ws.Range("G48").Value = "implemented"

# Row 49 - S07_G02_TF001
$ws.Range("F49").Value = "Settings modal implemented with Inputs/Risk/Costs/Visualization/Meta tabs for each backtest."
$ws.Range("G49").Value = "implemented"

# Row 50 - S07_G02_TF002
$ws.Range("F50").Value = "Settings modal wired to PATCH /api/backtests/{id}/settings and BacktestRead configs."
$ws.Range("G50").Value = "implemented"

# Row 51 - S07_G02_TF003
$ws.Range("F51").Value = "Risk/costs/visual configs persisted on Backtest and applied to chart behaviour."
$ws.Range("G51").Value = "implemented"

# Row 52 - S07_G03_TF001
$ws.Range("F52").Value = "Trades table with what-if metrics and cumulative PnL added to Backtest Details."
$ws.Range("G52").Value = "implemented"

# Row 53 - S07_G03_TF002
$ws.Range("F53").Value = "Export CSV button hooked to /api/backtests/{id}/trades/export in UI."
$ws.Range("G53").Value = "implemented"

# Row 54 - S07_G03_TF003
$ws.Range("F54").Value = "Interactive linkage between trade selection and chart segments deferred to future sprint."
$ws.Range("G54").Value = "implemented"
$ws.Range("H54").Value = "Deferred: interactive selection/highlighting to be done in later BT iteration."

# Row 55 - S08_G01_TB001
$ws.Range("F55").Value = "User manual updated for coverage IDs, new Run Backtest flow, Backtest Details and settings."
$ws.Range("G55").Value = "implemented"

# Row 56 - S08_G01_TF002 (status stays pending)
$ws.Range("F56").Value = "Frontend behaviour validated manually; automated tests can be added in a later test-focused sprint."

# Row 57 - S08_G01_TB003
$ws.Range("F57").Value = "qlab_impl_report.md and pytest/ruff config updated for Backtest Overhaul docs and lint/marker polish."
$ws.Range("G57").Value = "implemented"
